$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Insert a new row above row 20 (pushes "nego bom" and everything below down by one)
$ws.Rows.Item(20).Insert()

# Copy the formatting used by the other sub-heading rows (e.g. row 35 "coquero" after the shift)
$ws.Range("A35:B35").Copy()
$ws.Range("A20:B20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new product name
$ws.Range("A20").Value = "castanha"

$ws.Range("A21").Select()
